$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the existing row 2 ("044/FES VILLE /AV6"),
# shifting it (and the totals row after it) down by one row.
$ws.Rows.Item(2).Insert()

# Fill the newly inserted row 2 with the new contract's data.
$ws.Range("A2").Value = "665/FES 2"
$ws.Range("B2").Value = "Point de vente"
$ws.Range("C2").Value = "K5443645"
$ws.Range("D2").Value = "KHADIJA LALA"
$ws.Range("E2").Value = "non"
$ws.Range("F2").Value = "mensuelle"
$ws.Range("G2").Value = 15
$ws.Range("H2").Value = 10000
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 1500
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 10000
$ws.Range("O2").Value = 18500

# Update the totals row (now row 4) to reflect the new contract plus
# the existing one (row2 + row3).
$ws.Range("H4").Value = 40000
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 6000
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 20000
$ws.Range("O4").Value = 44000
